$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Muscles(L)")

# Row 2 (IP): add reference marks B..D
$ws.Range("B2").Value = 146
$ws.Range("C2").Value = 140
$ws.Range("D2").Value = 129

# Row 8 (BF): add reference marks B..G
$ws.Range("B8").Value = 255
$ws.Range("C8").Value = 260
$ws.Range("D8").Value = 263
$ws.Range("E8").Value = 258
$ws.Range("F8").Value = 257
$ws.Range("G8").Value = 262

# Row 11 (SMT): add reference marks B..D
$ws.Range("B11").Value = 264
$ws.Range("C11").Value = 255
$ws.Range("D11").Value = 248

# Row 14 (RF): add reference mark C
$ws.Range("C14").Value = 206

# Row 17 (VL): add reference marks C..E
$ws.Range("C17").Value = 153
$ws.Range("D17").Value = 152
$ws.Range("E17").Value = 175

# Row 26 (GN): update B and add reference marks C..E
$ws.Range("B26").Value = 214
$ws.Range("C26").Value = 220
$ws.Range("D26").Value = 214
$ws.Range("E26").Value = 220

# Row 29 (CT): add reference marks C..D
$ws.Range("C29").Value = 227
$ws.Range("D29").Value = 227

# Move the active tab / selection to the Muscles(L) sheet at T19
$ws.Activate()
$ws.Range("T19").Select()
